# Add 2022-Q4 data
#
# 1) Insert a new worksheet named "2022-Q4" right after "总计", populate it
#    with the new quarter's per-fund holdings data.
# 2) Insert a new row at the top of the "总计" (summary) sheet's data area
#    (row 2) with the 2022-Q4 summary figures, shifting existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q4" worksheet, positioned right after "总计".
# ---------------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $summarySheet)
$newSheet.Name = "2022-Q4"

# Match the outline settings (summaryBelow/summaryRight) used by the other
# sheets in this workbook.
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# Copy header formatting (bold / centered / bordered) from an existing
# quarter sheet so the new sheet matches the established look. (Column A
# has no header cell in the template, so copy B1:H1 only.)
$templateSheet = $wb.Worksheets.Item("2022-Q3")
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

# Header row text.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Copy the numeric-cell formatting (style index used by A2 in the other
# quarter sheets) onto A2 of the new sheet.
$templateSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)      # xlPasteFormats
$newSheet.Range("A2").Value = 0

# Data row (row 2) values. The text-like numeric columns (B, D-G) are
# stored as text in the source workbook (e.g. fund code "001513" has a
# leading zero, "基金规模" etc. are decimal-looking strings), so force a
# text number format before assignment to avoid Excel auto-converting
# them to numbers, then clear the format again so no extra style index
# is left behind on the cells (matching the unstyled look used elsewhere).
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "001513"
$newSheet.Range("C2").Value = "易方达信息产业混合"
$newSheet.Range("D2").Value = "33.11"
$newSheet.Range("E2").Value = "90.45"
$newSheet.Range("F2").Value = "1.83"
$newSheet.Range("G2").Value = "0.6059"
$newSheet.Range("B2:G2").ClearFormats()
$newSheet.Range("H2").Value = 8

# ---------------------------------------------------------------------------
# 2. Insert the 2022-Q4 summary row into the "总计" sheet.
# ---------------------------------------------------------------------------
$summarySheet.Rows(2).Insert()

# The insert operation can carry over header formatting into the new row;
# clear that so the B:D cells go back to the default (unstyled) look used
# by the rest of the data rows.
$summarySheet.Range("B2:D2").ClearFormats()

# Re-apply the "index" column style (used by A3:A6) onto the new A2 cell.
$summarySheet.Range("A3").Copy()
$summarySheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$summarySheet.Range("A2").Value = 0
$summarySheet.Range("B2").Value = "2022-Q4"
$summarySheet.Range("C2").Value = 1
$summarySheet.Range("D2").Value = 0.61

# The row-insert shifted the pre-existing index values in column A down
# along with the rest of the row content, but in the source data column A
# is simply the (row - 2) index, so renumber A3:A6 back to 1,2,3,4.
$summarySheet.Range("A3").Value = 1
$summarySheet.Range("A4").Value = 2
$summarySheet.Range("A5").Value = 3
$summarySheet.Range("A6").Value = 4

$excel.CutCopyMode = 0
